# create test case for all
#
# Adds email/password/firstName/lastName/phone/accountType columns (U:Z)
# to the "ALL" sheet, mirroring the sign-up test data already present on
# the "DataSignUp" sheet, and updates the remembered cell selection on
# three sheets (DataCustomer, ALL, DataSignUp).

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("ALL")

# --- seed U1:Z5 with the formatting already used by the neighbouring
#     Q:T columns (font/number-format), then overwrite the values. Using
#     two 3-column blocks (R:T, which style consistently for all 5 rows)
#     keeps the source/destination range sizes identical so the paste
#     doesn't smear across extra columns.
[void]$ws3.Range("R1:T5").Copy()
[void]$ws3.Range("U1:W5").PasteSpecial(-4122)
[void]$ws3.Range("R1:T5").Copy()
[void]$ws3.Range("X1:Z5").PasteSpecial(-4122)

# Header row
$ws3.Range("U1").Value = "email"
$ws3.Range("V1").Value = "password"
$ws3.Range("W1").Value = "firstName"
$ws3.Range("X1").Value = "lastName"
$ws3.Range("Y1").Value = "phone"
$ws3.Range("Z1").Value = "accountType"

# Row 2
$ws3.Range("U2").Value = "rinarcus@gmail.com"
$ws3.Range("V2").Value = "iniPassword"
$ws3.Range("W2").Value = "rin"
$ws3.Range("X2").Value = "arcus"
$ws3.Range("Y2").Value = 8798132
$ws3.Range("Z2").Value = "Customer"

# Row 3
$ws3.Range("U3").Value = "rinarcus@gmail.com"
$ws3.Range("V3").Value = "iniPassword"
$ws3.Range("W3").Value = "ris"
$ws3.Range("X3").Value = "kuy"
$ws3.Range("Y3").Value = 2121121
$ws3.Range("Z3").Value = "Customer"

# Row 4
$ws3.Range("U4").Value = "rinarcus@gmail.com"
$ws3.Range("V4").Value = "iniPassword"
$ws3.Range("W4").Value = "pan"
$ws3.Range("X4").Value = "day"
$ws3.Range("Y4").Value = 12122121
$ws3.Range("Z4").Value = "Customer"

# Row 5
$ws3.Range("U5").Value = "rinarcus@gmail.com"
$ws3.Range("V5").Value = "iniPassword"
$ws3.Range("W5").Value = "mae"
$ws3.Range("X5").Value = "mun"
$ws3.Range("Y5").Value = 11212111
$ws3.Range("Z5").Value = "Customer"

# --- update the remembered selections. Order matters: whichever sheet is
#     selected last becomes the saved active tab, so touch "ALL" last to
#     keep it the active sheet (as it was before the edit).
$ws1 = $wb.Worksheets.Item("DataCustomer")
[void]$ws1.Range("L1:L5").Select()

$ws5 = $wb.Worksheets.Item("DataSignUp")
[void]$ws5.Range("F1:F5").Select()

[void]$ws3.Range("L18").Select()
